# Fix "CO-CURRIULAR ACTIVITIES" -> "CO-CURRICULAR ACTIVITIES"
#
# The canonical diff shows the header run being re-typed as three
# adjacent runs: "CO-CURRI" + "C" + "ULAR ACTIVITIES" (i.e. a single
# "C" was typed in the middle of the existing word, which is exactly
# what happens when a human clicks between the "I" and the "U" and
# types a letter). None of the three runs carry any explicit run
# formatting, so we insert the missing "C" with change tracking on
# (which keeps Word from silently re-coalescing the edit back into a
# single run) and then accept just that one revision so no <w:ins>
# markup - or stray formatting - is left behind.

$d = $word.ActiveDocument

$oldText = "CO-CURRIULAR ACTIVITIES"
$splitAt = "CO-CURRI".Length   # insert the missing "C" right after "CO-CURRI"

# Find the paragraph that still has the typo (search instead of a
# hard-coded index so the script keeps working even if the document
# changes shape elsewhere).
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.Contains($oldText)) {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    $paraStart = $target.Range.Start
    $offset = $target.Range.Text.IndexOf($oldText)
    $insertPos = $paraStart + $offset + $splitAt

    $wasTracking = $d.TrackRevisions
    $d.TrackRevisions = $true

    $insPoint = $d.Range($insertPos, $insertPos)
    $insPoint.InsertAfter("C")

    # Accept only the single revision we just made (not AcceptAll, which
    # touches unrelated parts of the document) so the run split sticks
    # around without leaving <w:ins> markup behind.
    $d.Revisions(1).Accept()

    $d.TrackRevisions = $wasTracking
}
